$wb = $excel.ActiveWorkbook

# Sheet "f__Nocardioidaceae-b-p": remove the duplicated row 4
# (label_GCF_000312105_3.fasta) so the trailing duplicate rows 5-8
# shift up and become rows 4-7.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(4).Delete()

# Sheet "f__Propionibacteriaceae-b-p": remove the duplicated block of
# rows 14-26 so the remaining rows (27-81) shift up to become rows
# 14-68.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A14:E26").EntireRow.Delete()
